# New crime data collected — update the weekly CompStat 115th Precinct sheet:
#  - bump the report volume/number and the covered week dates
#  - refresh the Week-to-Date / 28-Day / Year-to-Date crime counts and the
#    derived percent-change columns for rows 15-30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text (shared strings) ----------------------------------------
# "Volume 30   Number  42" -> "... 43"
$ws.Range("A8").Value = "Volume 30   Number  43"

# "Report Covering the Week  10/16/2023  Through  10/22/2023" -> next week
$ws.Range("C9").Value = "Report Covering the Week  10/23/2023  Through  10/29/2023"

# ---- Row 15: Rape ----------------------------------------------------------
$ws.Range("C15").Value = "'0"
$ws.Range("A15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("I15").Value = 21
$ws.Range("J15").Value = 35
$ws.Range("K15").Value = -40
$ws.Range("L15").Value = -16
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -22.222222222222

# ---- Row 16: Robbery -------------------------------------------------------
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 12
$ws.Range("F16").Value = 43
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = 34.375
$ws.Range("I16").Value = 295
$ws.Range("J16").Value = 249
$ws.Range("K16").Value = 18.473895582329
$ws.Range("L16").Value = 48.989898989899
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -72.811059907834

# ---- Row 17: Fel. Assault ---------------------------------------------------
$ws.Range("C17").Value = 9
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = -3.125
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 355
$ws.Range("K17").Value = 12.676056338028
$ws.Range("L17").Value = 24.223602484472
$ws.Range("M17").Value = 47.058823529411
$ws.Range("N17").Value = 9.890109890109

# ---- Row 18: Burglary -------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -45.454545454545
$ws.Range("I18").Value = 120
$ws.Range("J18").Value = 125
$ws.Range("K18").Value = -4
$ws.Range("L18").Value = 2.564102564102
$ws.Range("M18").Value = -51.807228915662
$ws.Range("N18").Value = -92.757996378998

# ---- Row 19: Gr. Larceny ----------------------------------------------------
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = 13.636363636363
$ws.Range("F19").Value = 83
$ws.Range("G19").Value = 88
$ws.Range("H19").Value = -5.681818181818
$ws.Range("I19").Value = 708
$ws.Range("J19").Value = 815
$ws.Range("K19").Value = -13.128834355828
$ws.Range("L19").Value = 58.035714285714
$ws.Range("M19").Value = 71.428571428571
$ws.Range("N19").Value = -41.487603305785

# ---- Row 20: G.L.A. ----------------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -28.571428571428
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 12
$ws.Range("I20").Value = 284
$ws.Range("J20").Value = 260
$ws.Range("K20").Value = 9.230769230769
$ws.Range("L20").Value = 64.161849710982
$ws.Range("M20").Value = 51.063829787234
$ws.Range("N20").Value = -84.925690021231

# ---- Row 21: TOTAL ------------------------------------------------------------
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = -5.454545454545
$ws.Range("F21").Value = 195
$ws.Range("G21").Value = 193
$ws.Range("H21").Value = 1.036269430051
$ws.Range("I21").Value = 1830
$ws.Range("J21").Value = 1842
$ws.Range("K21").Value = -0.651465798045
$ws.Range("L21").Value = 42.191142191142
$ws.Range("M21").Value = 26.556016597510
$ws.Range("N21").Value = -70.682473566164

# ---- Row 22: Transit ------------------------------------------------------------
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 7
$ws.Range("I22").Value = 74
$ws.Range("J22").Value = 54
$ws.Range("K22").Value = 37.037037037037
$ws.Range("L22").Value = 208.333333333333
$ws.Range("M22").Value = 174.074074074074

# Row 23 (Housing) is unchanged.

# ---- Row 24: Petit Larceny -------------------------------------------------------
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -13.888888888888
$ws.Range("F24").Value = 148
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = -11.377245508982
$ws.Range("I24").Value = 1641
$ws.Range("J24").Value = 1567
$ws.Range("K24").Value = 4.722399489470
$ws.Range("L24").Value = 45.866666666666
$ws.Range("M24").Value = 74.019088016967

# ---- Row 25: Misd. Assault -------------------------------------------------------
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = 31.578947368421
$ws.Range("F25").Value = 76
$ws.Range("G25").Value = 88
$ws.Range("H25").Value = -13.636363636363
$ws.Range("I25").Value = 808
$ws.Range("J25").Value = 751
$ws.Range("K25").Value = 7.589880159786
$ws.Range("L25").Value = 15.099715099715
$ws.Range("M25").Value = 2.668360864040

# ---- Row 26: UCR Rape* -------------------------------------------------------
$ws.Range("C26").Value = "'0"
$ws.Range("A26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 41
$ws.Range("J26").Value = 46
$ws.Range("K26").Value = -10.869565217391
$ws.Range("L26").Value = -6.818181818181

# ---- Row 27: Other Sex Crimes -------------------------------------------------------
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 200
$ws.Range("I27").Value = 128
$ws.Range("J27").Value = 86
$ws.Range("K27").Value = 48.837209302325
$ws.Range("L27").Value = 42.222222222222

# ---- Row 28: Shooting Vic. -------------------------------------------------------
$ws.Range("N28").Value = -90.566037735849

# ---- Row 29: Shooting Inc. -------------------------------------------------------
$ws.Range("N29").Value = -91.836734693877

# ---- Row 30: Hate Crimes -------------------------------------------------------
$ws.Range("F30").Value = 1
